$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: swap D3 and E3 team values (G2 <-> Fnatic)
$ws.Range("D3").Value = "Fnatic"
$ws.Range("E3").Value = "G2"

# Row 5: update score values
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 2

# Row 6: update score values
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 1

# Update the active selection to B8, matching the saved sheet view state
$ws.Range("B8").Select()
